# Adds the "Branches" section (rows 4-12) to Sheet1, per commit:
# "Commands of Branch entered please check"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlHAlignLeft = -4131
$xlHAlignCenter = -4108
$xlVAlignCenter = -4108

function Set-CellStyle {
    param($range, $bold, $hAlign, $wrap)
    $range.Font.Name = "Times New Roman"
    $range.Font.Size = 11
    $range.Font.Bold = $bold
    $range.HorizontalAlignment = $hAlign
    $range.VerticalAlignment = $xlVAlignCenter
    $range.WrapText = $wrap
}

# --- Enter values in the same order the original author typed them ---
# (this keeps the shared-string table / cellXf allocation order aligned)

# Row 4: "Branches" section title (bold, left aligned)
$ws.Range("A4").Value = "Branches"

# Row 5
$ws.Range("A5").Value = "To create a new branch"
$ws.Range("B5").Value = "Branch command to be used"
$ws.Range("C5").Value = "git branch 'branch_name'"

# Row 6 (command columns typed first, description last)
$ws.Range("B6").Value = "switch/checkout"
$ws.Range("C6").Value = "git switch 'branch_name'/ git checkout 'branch_name'"
$ws.Range("A6").Value = "To change to other branch"

# Row 7 / Row 8 (column A for both rows, then C7, then back for B7/B8/C8)
$ws.Range("A7").Value = "To create a new branch and change to it at once"
$ws.Range("C7").Value = "1. git switch -c 'branch_name' `n2.git checkout -b 'branch_name'"
$ws.Range("A8").Value = "To list branches that are currently in working mode"
$ws.Range("B7").Value = "can be done either using switch -c and checkout -b"
$ws.Range("B8").Value = "use branch"
$ws.Range("C8").Value = "git branch "

# Row 9
$ws.Range("A9").Value = "To list all the branches"
$ws.Range("B9").Value = "use branch -a"
$ws.Range("C9").Value = "git branch -a"

# Row 10
$ws.Range("A10").Value = "To rename the branches"
$ws.Range("B10").Value = "branch -m"
$ws.Range("C10").Value = "git branch -m 'old branch name' 'new branch name'"

# Row 11
$ws.Range("A11").Value = "To delete a particular branch"
$ws.Range("B11").Value = "branch -d"
$ws.Range("C11").Value = "git branch -d 'branch name'"

# Row 12
$ws.Range("A12").Value = "To do forceful deletion "
$ws.Range("B12").Value = "branch -D"
$ws.Range("C12").Value = "git branch -D 'branch name'"

# --- Apply formatting ---

Set-CellStyle $ws.Range("A4") $true $xlHAlignLeft $false

Set-CellStyle $ws.Range("A5") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("B5") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("C5") $true $xlHAlignCenter $false

Set-CellStyle $ws.Range("A6") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("B6") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("C6") $true $xlHAlignCenter $false

Set-CellStyle $ws.Range("A7") $false $xlHAlignLeft $true
Set-CellStyle $ws.Range("B7") $false $xlHAlignLeft $true
Set-CellStyle $ws.Range("C7") $true $xlHAlignLeft $true
$ws.Rows.Item(7).RowHeight = 27.6

Set-CellStyle $ws.Range("A8") $false $xlHAlignLeft $true
Set-CellStyle $ws.Range("B8") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("C8") $true $xlHAlignLeft $false
$ws.Rows.Item(8).RowHeight = 27.6

Set-CellStyle $ws.Range("A9") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("B9") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("C9") $true $xlHAlignLeft $false

Set-CellStyle $ws.Range("A10") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("B10") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("C10") $true $xlHAlignLeft $false

Set-CellStyle $ws.Range("A11") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("B11") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("C11") $false $xlHAlignLeft $false

Set-CellStyle $ws.Range("A12") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("B12") $false $xlHAlignLeft $false
Set-CellStyle $ws.Range("C12") $false $xlHAlignLeft $false

# Final selection matches the author's last position (B12)
$ws.Range("B12").Select()
